# Insert a new column before column C and add the "tipo.dado" header,
# shifting the existing pub.pre / pub.ead / priv.pre / priv.ead columns
# one position to the right (C->D, D->E, E->F, F->G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C").Insert()
$ws.Range("C1").Value = "tipo.dado"
